# Daily attendance processing - 2025-12-31 21:54:49
# Normalizes the "Recorded By" (column G) lists so that the first
# recorder in the comma-separated list is rotated to the end, unless
# the list already ends with the exact (case-sensitive) token "System".

function Test-ExactEquals($s1, $s2) {
    if ($s1.Length -ne $s2.Length) {
        return $false
    }
    for ($i = 0; $i -lt $s1.Length; $i++) {
        $c1 = [int][char]$s1[$i]
        $c2 = [int][char]$s2[$i]
        if ($c1 -ne $c2) {
            return $false
        }
    }
    return $true
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -eq $null) {
        continue
    }

    if ($value -like "*,*") {
        $parts = $value -split ", "
        if ($parts.Count -gt 1) {
            $firstPart = $parts[0]
            $lastPart = $parts[$parts.Count - 1]

            if (-not (Test-ExactEquals $lastPart "System")) {
                $rest = $parts[1..($parts.Count - 1)]
                $newParts = $rest + @($firstPart)
                $newValue = $newParts -join ", "
                $cell.Value = $newValue
            }
        }
    }
}
